$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4218206666666666
$ws.Range("H2").Value = 1.265462
$ws.Range("I2").Value = 0.204479520571209
$ws.Range("J2").Value = 0.204479520571209
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.312815
$ws.Range("N2").Value = 0.938445
$ws.Range("O2").Value = 0.008213170494906699
$ws.Range("P2").Value = 0.008213170494906699
$ws.Range("Q2").Value = 0.1319518318433333
$ws.Range("R2").Value = 1.18756648659
$ws.Range("S2").Value = 0.001679425165168121
$ws.Range("T2").Value = 0.001679425165168121

$ws.Range("G3").Value = 0.4218206666666666
$ws.Range("H3").Value = 1.265462
$ws.Range("I3").Value = 0.204479520571209
$ws.Range("J3").Value = 0.204479520571209
$ws.Range("O3").Value = 0.6524076620340182
$ws.Range("P3").Value = 0.6524076620340182
$ws.Range("Q3").Value = 10.48150481807244
$ws.Range("R3").Value = 94.33354336265199
$ws.Range("S3").Value = 0.1334040059496994
$ws.Range("T3").Value = 0.1334040059496994

$ws.Range("G4").Value = 0.4218206666666666
$ws.Range("H4").Value = 1.265462
$ws.Range("I4").Value = 0.204479520571209
$ws.Range("J4").Value = 0.204479520571209
$ws.Range("M4").Value = 12.866992
$ws.Range("N4").Value = 38.600976
$ws.Range("O4").Value = 0.3378316226926476
$ws.Range("P4").Value = 0.3378316226926476
$ws.Range("Q4").Value = 5.427563143434667
$ws.Range("R4").Value = 48.848068290912
$ws.Range("S4").Value = 0.06907964824198615
$ws.Range("T4").Value = 0.06907964824198615

$ws.Range("G5").Value = 0.4218206666666666
$ws.Range("H5").Value = 1.265462
$ws.Range("I5").Value = 0.204479520571209
$ws.Range("J5").Value = 0.204479520571209
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05894133333333334
$ws.Range("N5").Value = 0.176824
$ws.Range("O5").Value = 0.001547544778427486
$ws.Range("P5").Value = 0.001547544778427486
$ws.Range("Q5").Value = 0.02486267252088889
$ws.Range("R5").Value = 0.223764052688
$ws.Range("S5").Value = 0.0003164412143553302
$ws.Range("T5").Value = 0.0003164412143553302

$ws.Range("I6").Value = 0.3030684321645684
$ws.Range("J6").Value = 0.3030684321645683
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.312815
$ws.Range("N6").Value = 0.938445
$ws.Range("O6").Value = 0.008213170494906699
$ws.Range("P6").Value = 0.008213170494906699
$ws.Range("Q6").Value = 0.1955718337283333
$ws.Range("R6").Value = 1.760146503555
$ws.Range("S6").Value = 0.002489152704991665
$ws.Range("T6").Value = 0.002489152704991665

$ws.Range("I7").Value = 0.3030684321645684
$ws.Range("J7").Value = 0.3030684321645683
$ws.Range("O7").Value = 0.6524076620340182
$ws.Range("P7").Value = 0.6524076620340182
$ws.Range("S7").Value = 0.1977241672648015
$ws.Range("T7").Value = 0.1977241672648015

$ws.Range("I8").Value = 0.3030684321645684
$ws.Range("J8").Value = 0.3030684321645683
$ws.Range("M8").Value = 12.866992
$ws.Range("N8").Value = 38.600976
$ws.Range("O8").Value = 0.3378316226926476
$ws.Range("P8").Value = 0.3378316226926476
$ws.Range("Q8").Value = 8.044439109402667
$ws.Range("R8").Value = 72.399951984624
$ws.Range("S8").Value = 0.1023861002250727
$ws.Range("T8").Value = 0.1023861002250727

$ws.Range("I9").Value = 0.3030684321645684
$ws.Range("J9").Value = 0.3030684321645683
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.05894133333333334
$ws.Range("N9").Value = 0.176824
$ws.Range("O9").Value = 0.001547544778427486
$ws.Range("P9").Value = 0.001547544778427486
$ws.Range("Q9").Value = 0.03685010195288889
$ws.Range("R9").Value = 0.331650917576
$ws.Range("S9").Value = 0.0004690119697024826
$ws.Range("T9").Value = 0.0004690119697024825

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.469433
$ws.Range("H10").Value = 1.408299
$ws.Range("I10").Value = 0.2275598195290835
$ws.Range("J10").Value = 0.2275598195290835
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.312815
$ws.Range("N10").Value = 0.938445
$ws.Range("O10").Value = 0.008213170494906699
$ws.Range("P10").Value = 0.008213170494906699
$ws.Range("Q10").Value = 0.146845683895
$ws.Range("R10").Value = 1.321611155055
$ws.Range("S10").Value = 0.001868987595582562
$ws.Range("T10").Value = 0.001868987595582562

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.469433
$ws.Range("H11").Value = 1.408299
$ws.Range("I11").Value = 0.2275598195290835
$ws.Range("J11").Value = 0.2275598195290835
$ws.Range("O11").Value = 0.6524076620340182
$ws.Range("P11").Value = 0.6524076620340182
$ws.Range("Q11").Value = 11.66458791633933
$ws.Range("R11").Value = 104.981291247054
$ws.Range("S11").Value = 0.1484617698318525
$ws.Range("T11").Value = 0.1484617698318525

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.469433
$ws.Range("H12").Value = 1.408299
$ws.Range("I12").Value = 0.2275598195290835
$ws.Range("J12").Value = 0.2275598195290835
$ws.Range("M12").Value = 12.866992
$ws.Range("N12").Value = 38.600976
$ws.Range("O12").Value = 0.3378316226926476
$ws.Range("P12").Value = 0.3378316226926476
$ws.Range("Q12").Value = 6.040190655536001
$ws.Range("R12").Value = 54.36171589982401
$ws.Range("S12").Value = 0.07687690309115633
$ws.Range("T12").Value = 0.0768769030911563

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.469433
$ws.Range("H13").Value = 1.408299
$ws.Range("I13").Value = 0.2275598195290835
$ws.Range("J13").Value = 0.2275598195290835
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.05894133333333334
$ws.Range("N13").Value = 0.176824
$ws.Range("O13").Value = 0.001547544778427486
$ws.Range("P13").Value = 0.001547544778427486
$ws.Range("Q13").Value = 0.02766900693066667
$ws.Range("R13").Value = 0.249021062376
$ws.Range("S13").Value = 0.0003521590104921343
$ws.Range("T13").Value = 0.0003521590104921342

$ws.Range("G14").Value = 0.546446
$ws.Range("H14").Value = 1.639338
$ws.Range("I14").Value = 0.2648922277351392
$ws.Range("J14").Value = 0.2648922277351391
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.312815
$ws.Range("N14").Value = 0.938445
$ws.Range("O14").Value = 0.008213170494906699
$ws.Range("P14").Value = 0.008213170494906699
$ws.Range("Q14").Value = 0.17093650549
$ws.Range("R14").Value = 1.53842854941
$ws.Range("S14").Value = 0.002175605029164351
$ws.Range("T14").Value = 0.00217560502916435

$ws.Range("G15").Value = 0.546446
$ws.Range("H15").Value = 1.639338
$ws.Range("I15").Value = 0.2648922277351392
$ws.Range("J15").Value = 0.2648922277351391
$ws.Range("O15").Value = 0.6524076620340182
$ws.Range("P15").Value = 0.6524076620340182
$ws.Range("Q15").Value = 13.57822609090533
$ws.Range("R15").Value = 122.204034818148
$ws.Range("S15").Value = 0.1728177189876648
$ws.Range("T15").Value = 0.1728177189876648

$ws.Range("G16").Value = 0.546446
$ws.Range("H16").Value = 1.639338
$ws.Range("I16").Value = 0.2648922277351392
$ws.Range("J16").Value = 0.2648922277351391
$ws.Range("M16").Value = 12.866992
$ws.Range("N16").Value = 38.600976
$ws.Range("O16").Value = 0.3378316226926476
$ws.Range("P16").Value = 0.3378316226926476
$ws.Range("Q16").Value = 7.031116310432001
$ws.Range("R16").Value = 63.280046793888
$ws.Range("S16").Value = 0.08948897113443241
$ws.Range("T16").Value = 0.0894889711344324

$ws.Range("G17").Value = 0.546446
$ws.Range("H17").Value = 1.639338
$ws.Range("I17").Value = 0.2648922277351392
$ws.Range("J17").Value = 0.2648922277351391
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.05894133333333334
$ws.Range("N17").Value = 0.176824
$ws.Range("O17").Value = 0.001547544778427486
$ws.Range("P17").Value = 0.001547544778427486
$ws.Range("Q17").Value = 0.03220825583466667
$ws.Range("R17").Value = 0.289874302512
$ws.Range("S17").Value = 0.0004099325838775391
$ws.Range("T17").Value = 0.000409932583877539
